$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed crypto market data (Price + Volume(1h) change%),
# and swap the Toncoin/Dogecoin and Bittensor/Maker row ordering to match
# the latest coinranking.com pull.
#
# Column D (Price) values are written with a leading apostrophe (escaped
# with a backtick in the double-quoted string) so numeric-looking text
# (e.g. "1.00", "47.10", or thousand-grouped "67.194.11") is stored
# verbatim as text instead of being parsed/normalized into a number.
$ws.Range('D2').Value = "`'67.194.11"
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = "`'3.134.72"
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('D4').Value = "`'1.00"
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "`'581.43"
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Value = "`'174.18"
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('D7').Value = "`'1.00"
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "`'0.156"
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').Value = "`'6.44"
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').Value = "`'0.481"
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = "`'0.0000250"
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = "`'37.73"
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').Value = "`'0.122"
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = "`'67.120.29"
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = "`'7.16"
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').Value = "`'3.132.88"
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').Value = "`'16.40"
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('D19').Value = "`'492.13"
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').Value = "`'0.711"
$ws.Range('D21').Value = "`'7.89"
$ws.Range('E21').Value = '  +4.97%  '
$ws.Range('D22').Value = "`'84.34"
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').Value = "`'13.35"
$ws.Range('E23').Value = '  +2.67%  '
$ws.Range('E24').Value = '  -1.75%  '
$ws.Range('D25').Value = "`'10.39"
$ws.Range('E25').Value = '  +3.43%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = "`'7.96"
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('E29').Value = '  +0.46%  '
$ws.Range('D30').Value = "`'28.84"
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = "`'0.0₃0958"
$ws.Range('E32').Value = '  -5.31%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').Value = "`'0.981"
$ws.Range('E35').Value = '  -2.18%  '
$ws.Range('D36').Value = "`'47.10"
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  +1.85%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = "`'2.837.13"
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = "`'387.06"
$ws.Range('E43').Value = '  +2.12%  '
$ws.Range('D44').Value = "`'2.60"
$ws.Range('E44').Value = '  -6.95%  '
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').Value = "`'136.09"
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = "`'25.12"
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').Value = "`'6.81"
$ws.Range('E51').Value = '  -0.05%  '
